# Trade #120 closed at 2026-02-17 09:28:58 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up figures to account for the
# newly closed trade, and appends the new trade row to both the
# "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet roll-up numbers
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.68            # Current Capital
$summary.Range("B4").Value = 0.6899999999999999 # Total P&L $
$summary.Range("B5").Value = 0.11               # Total P&L %
$summary.Range("B6").Value = 120                # Total Trades
$summary.Range("B8").Value = 45                 # Losing Trades
$summary.Range("B9").Value = 45                 # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.68              # Capital
$status.Range("D4").Value = 120                 # Trades
$status.Range("E4").Value = 0.6899999999999999  # P&L $
$status.Range("F4").Value = 0.68                # P&L %
$status.Range("G4").Value = 45                  # Win Rate %

# ---------------------------------------------------------------------
# 3. Append new trade #120 row to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------
$newRow = 121

function Add-TradeRow($ws, $row) {
    $ws.Range("A$row").Value = 120
    # Leading apostrophe forces these to stay plain text instead of being
    # auto-parsed into Excel date/time serial numbers.
    $ws.Range("B$row").Value = "'2026-02-17"
    $ws.Range("C$row").Value = "'09:28:52"
    $ws.Range("D$row").Value = "MarketMaking"
    $ws.Range("E$row").Value = "DOWN"
    $ws.Range("F$row").Value = 0.824834
    $ws.Range("G$row").Value = 0.813923
    $ws.Range("H$row").Value = "CLOSED"
    $ws.Range("I$row").Value = -1.3227
    $ws.Range("J$row").Value = -0.01
    $ws.Range("K$row").Value = 100.68
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P$row").Value = "early_exit"
    $ws.Range("Q$row").Value = 0.16
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades $newRow

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking $newRow

Write-Host "Trade #120 appended; summary figures updated."
